# Apply the "Deploying to gh-pages" metadata refresh to the
# StructureDefinition-restricted-benefits workbook:
#   - Metadata sheet: URL / Version / Date / Publisher bump
#   - Elements sheet: clear the root Extension row's Constraint(s) cell

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/restricted-benefits"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
# Row 2 is the root "Extension" element; column AI is "Constraint(s)".
# It no longer carries the ele-1/ext-1 FHIRPath constraint text.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
